$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Replace the whole paragraph's XML (keeping its pPr and leading empty run
# intact) so that the <w:r/> placeholder run already present in the source
# paragraph is not merged away by a plain Find/Replace text substitution.
function Replace-BulletParagraph($oldText, $newText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$oldText*") {
            $pPrXml = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'
            $fragment = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                        '<w:p>' + $pPrXml + '<w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></pkg:xmlData>'
            $p.Range.InsertXML($fragment)
            return $true
        }
    }
    return $false
}

# Title (appears twice: Heading1 and the bold run near the end) - plain
# Find/Replace keeps the existing structure intact for both occurrences.
Replace-Text "Play Diwinity Slot for Free - Review and Bonus Features" "Play Diwinity Slot Game for Free"

# "What we like" bullet list
Replace-BulletParagraph "Low volatility with frequent wins" "Sleek graphics with cartoon-style symbols"
Replace-BulletParagraph "Two Wild symbols and a Wild Scatter" "Low volatility for frequent wins"
Replace-BulletParagraph "Bonus game with cash prizes" "Two bonus features - Free Spins and Bonus game"
Replace-BulletParagraph "Option to choose between Free Spins or the Bonus game" "Maximum win of 8x with x2 multipliers"

# "What we don't like" bullet list
Replace-BulletParagraph "RTP value of 94.88% is lower than some other slots" "RTP value is 94.88%"
Replace-BulletParagraph "Graphics may not appeal to everyone" "Limited number of paylines"

# Meta description (italic run)
Replace-Text "Read our review of Diwinity online slot game and learn about its key features and bonus features. Play for free and choose between Free Spins or the Bonus game." "Experience the ancient Greek gods in Diwinity, a free slot game with exciting bonus features."
